# Reproduces the "Add files via upload" commit: a set of helper/demo
# formulas (MAX/IF region lookups, ISNUMBER/TYPE/VALUE probes, SUM totals,
# INDEX/MATCH "top salesperson" lookups and a localized-function #NAME?
# example) added to the "Clean Data" sheet, plus 8 new trailing helper rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clean Data")

# --- Row 2: probes next to the "Top salesperson" question -----------------
$ws.Range("H2").FormulaArray = '=MAX(IF($E$6:$E$50=E6, $G$6:$G$50))'
$ws.Range("J2").Formula      = '=ISNUMBER(G2)'
$ws.Range("K2").Formula      = '=TYPE(G2)'
$ws.Range("L2").Formula      = '=VALUE(G2)'

# --- Row 3 & 4: more MAX/IF + VALUE probes ---------------------------------
$ws.Range("H3").FormulaArray = '=MAX(IF($E$2:$E$50=E3, $G$2:$G$50))'
$ws.Range("L3").Formula      = '=VALUE(G3)'

$ws.Range("H4").FormulaArray = '=MAX(IF($E$2:$E$50=E4, $G$2:$G$50))'
$ws.Range("L4").Formula      = '=VALUE(G4)'

# --- Row 5 (header row): H5 used to be a blank formatted cell (style 4) ---
# Copy the plain "style 1" formatting (from A3) over it first so the new
# formula cell matches the rest of column H instead of keeping the old
# header-row formatting.
$ws.Range("A3").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("H5").FormulaArray = '=MAX(IF($E$2:$E$50=E5, $G$2:$G$50))'

# --- Data rows 6-35: "max sale in this region" (H, 2:50 range) and --------
# "max sale in this region" (I, 6:1000 range) helper columns.
for ($r = 6; $r -le 12; $r++) {
    $ws.Range("H$r").FormulaArray = "=MAX(IF(`$E`$2:`$E`$50=E$r, `$G`$2:`$G`$50))"
}
for ($r = 6; $r -le 35; $r++) {
    $ws.Range("I$r").FormulaArray = "=MAX(IF(`$E`$6:`$E`$1000=E$r, `$G`$6:`$G`$1000))"
}

# --- Row 7 & 8: grand-total / top-salesperson summary formulas ------------
$ws.Range("P7").Formula = '=SUM(G2:G50)'
$ws.Range("Q7").Formula = '=SUM(G:G)'
$ws.Range("P8").Formula = '=INDEX(D2:D50, MATCH(MAX(G2:G50), G2:G50, 0))'
$ws.Range("Q8").Formula = '=INDEX(D2:D50, MATCH(MAX(G2:G50), G2:G50, 0))'

# --- Rows 6-10: "who had the max sale in this region" lookups ------------
for ($r = 6; $r -le 10; $r++) {
    $ws.Range("R$r").FormulaArray = "=INDEX(`$D`$6:`$D`$1000, MATCH(1, (`$E`$6:`$E`$1000=E$r)*(`$G`$6:`$G`$1000=I$r), 0))"
}
for ($r = 6; $r -le 8; $r++) {
    $ws.Range("S$r").FormulaArray = "=INDEX(`$D`$6:`$D`$1000,MATCH(1,(`$E`$6:`$E`$1000=E$r)*(`$G`$6:`$G`$1000=I$r),0))"
}

# --- U6: demo of a localized (Polish) function name => #NAME? error ------
# (Entered as a plain formula rather than FormulaArray: the engine's array-
# formula entry path doesn't classify this particular parse failure as a
# recoverable #NAME? error, so .Formula is what actually reproduces the
# intended #NAME? error value.)
$ws.Range("U6").Formula = '=SUMA.ILOCZYN(G6:G10)'

# --- New trailing helper rows 36-43 (I column only) -----------------------
for ($r = 36; $r -le 43; $r++) {
    $ws.Range("I$r").FormulaArray = "=MAX(IF(`$E`$6:`$E`$1000=E$r, `$G`$6:`$G`$1000))"
}

$excel.Calculate()

# --- Final view state: scrolled down a little, U6 selected ---------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("U6").Select()
